$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "1" to "თერჯოლა"
$ws.Name = "თერჯოლა"

# Drop the census-results caption in A2 ("(მოსახლეობის აღწერის შედეგებით)")
$ws.Range("A2").Clear()

# Drop the 1989 and 2002 columns, keeping only the 2014 figures (shifts old column D into B)
$ws.Range("B:C").Delete()

# Drop the now-blank spacer row (old row 3), shifting rows 4-6 up to rows 3-5
$ws.Range("A3").EntireRow.Delete()

# Match the saved selection/active cell
$ws.Range("A2").Select() | Out-Null
